# Upload Leave Card 12/27/2023 4:01 PM
# A new leave-card entry is inserted into the Table1 ListObject (Sheet1)
# as row 535, pushing every subsequent row down by one. The row above
# (534) and the two rows right below the newly inserted one pick up the
# actual leave entries (SL / FL) that were recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Insert a new physical row above the old row 535 - this shifts every
#    row from 535..711 down to 536..712, values/formulas intact.
$ws.Rows.Item(535).Insert()

# 2) The freshly inserted row has no formatting (default/no borders).
#    Row 533 is an existing "blank PERIOD" leave-entry row with exactly
#    the formatting this new row needs (A blank-date style, K date style,
#    etc.), so copy its formats down.
$ws.Range("A533:K533").Copy()
$ws.Range("A535:K535").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Grow Table1 so it covers the new row (A8:K712).
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K712"))

# 4) Restore the calculated-column formula in the new row's "EARNED "
#    column (PasteSpecial only carried formats, not formulas/values).
$ws.Range("G535").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),""."",Table1[[#This Row],[EARNED]])"
$ws.Range("G535").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# 5) Row 534 (10/01/2023 period) now records a Sick Leave usage:
#    SL(1-0-0), EARNED 1.25, Absence Undertime W/Pay(hrs)=1, dated 10/13/2023.
$ws.Range("B534").Value = "SL(1-0-0)"
$ws.Range("C534").Value = 1.25
$ws.Range("H534").Value = 1
$ws.Range("K534").Value = 45212

# 6) The newly inserted row 535 records a Forced Leave: FL(1-0-0),
#    1 day undertime, dated 10/26/2023.
$ws.Range("B535").Value = "FL(1-0-0)"
$ws.Range("D535").Value = 1
$ws.Range("K535").Value = 45225

# 7) Row 536 (formerly row 535, shifted down) also records FL(1-0-0),
#    1 day undertime, dated 11/23/2023.
$ws.Range("B536").Value = "FL(1-0-0)"
$ws.Range("D536").Value = 1
$ws.Range("K536").Value = 45253

# 8) Row 537 (formerly row 536, shifted down) records FL(3-0-0),
#    3 days undertime, with a text remark instead of a date.
$ws.Range("B537").Value = "FL(3-0-0)"
$ws.Range("D537").Value = 3
$ws.Range("K537").Value = "12/27-29/2023"

# 9) Move the selection to mirror where the editor ended up (sheet view
#    pane/selection tracked in the diff).
$ws.Range("F537").Select()
